$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 28 (Dur. Order / DGORDER - M/M) ---
$ws.Range("F28").Value = 0.02880223870746335
$ws.Range("G28").Value = -0.02799901206372835

# --- Row 29 (DGORDER - Y/Y  |  T5YIFR) ---
$ws.Range("F29").Value = 0.07530887191904141
$ws.Range("G29").Value = 0.03341358778313566
$ws.Range("N29").Value = 45979
$ws.Range("Q29").Value = 2.18
$ws.Range("R29").Value = 2.19
$ws.Range("S29").Value = $null
$ws.Range("T29").Value = $null
$ws.Range("U29").Value = 2.18

# --- Row 30 (ADXDNO - M/M  |  T10YIE) ---
$ws.Range("F30").Value = 0.01912935471760346
$ws.Range("G30").Value = -0.02404555711932721
$ws.Range("N30").Value = 45979
$ws.Range("Q30").Value = 2.27
$ws.Range("R30").Value = 2.28
$ws.Range("S30").Value = $null
$ws.Range("T30").Value = $null
$ws.Range("U30").Value = 2.28

# --- Row 31 (ADXDNO - Y/Y) ---
$ws.Range("F31").Value = 0.06676582929337446
$ws.Range("G31").Value = 0.0329297153895499

# Highlight the four "as-of" dates in column C (rows 28-31) the same way
# column N's updated dates are highlighted - copy N29's format (style 48,
# yellow fill) over C28:C31 without touching their values.
$ws.Range("N29").Copy()
$ws.Range("C28:C31").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 47 (FFR / DFF) ---
$ws.Range("N47").Value = 45978
$ws.Range("S47").Value = 3.88
$ws.Range("T47").Value = 3.88
$ws.Range("U47").Value = 3.88

# --- Row 48 (2y UST / DGS2) ---
$ws.Range("N48").Value = 45978
$ws.Range("Q48").Value = 3.6
$ws.Range("R48").Value = $null
$ws.Range("S48").Value = $null
$ws.Range("T48").Value = 3.62

# --- Row 49 (5y UST / DGS5) ---
$ws.Range("N49").Value = 45978
$ws.Range("Q49").Value = 3.72
$ws.Range("R49").Value = $null
$ws.Range("S49").Value = $null
$ws.Range("T49").Value = 3.74
$ws.Range("U49").Value = 3.71

# --- Row 50 (10y UST / DGS10) ---
$ws.Range("N50").Value = 45978
$ws.Range("Q50").Value = 4.13
$ws.Range("R50").Value = $null
$ws.Range("S50").Value = $null
$ws.Range("T50").Value = 4.14
$ws.Range("U50").Value = 4.11

# --- Row 52 (BAA / DBAA) ---
$ws.Range("N52").Value = 45978
$ws.Range("Q52").Value = 5.9
$ws.Range("R52").Value = $null
$ws.Range("S52").Value = $null
$ws.Range("T52").Value = 5.91
$ws.Range("U52").Value = 5.88
